$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K5").Value = 17.11
$ws.Range("L5").Value = 24.98

$ws.Range("K7").Value = 90.47799999999999
$ws.Range("L7").Value = 68.64

$ws.Range("K8").Value = 145.588
$ws.Range("L8").Value = 181.74
